$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the TOLL/LANES values that were swapped between B10 and B11 ---
$ws.Range("B10").Value = "LANES"
$ws.Range("B11").Value = "TOLL"

# --- Add a new "FTC2" column (column D) alongside the existing Network/CPP
#     Fields table (rows 2-10), flagging each listed field with a 1. The
#     FTYPE row (row 6) gets the literal label "FTC2" instead of a flag. ---

# Rows 2-5: plain "1"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1

# Row 6: FTC2 label
$ws.Range("D6").Value = "FTC2"

# Rows 7-9: "1"
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1

# Row 10: plain "1"
$ws.Range("D10").Value = 1

# --- Move the active selection to B10, matching the editor's last position ---
$ws.Range("B10").Select()
